# Updated PCM playback metrics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the D (Theoretical us/interrupt offset) and F (Real us/interrupt)
#     formulas for every data row: the constant changes from 66 to 48.
for ($r = 3; $r -le 14; $r++) {
    $ws.Range("D$r").Formula = "=C$r-48"
    $ws.Range("F$r").Formula = "=ROUND(A$r/(D$r+48),0)"
}

# --- Notes column (G): a new note "Maximum, with an empty stage" is
#     inserted before "Maximum, without animations.", pushing every
#     subsequent note down one row, and a brand new trailing note
#     "Sounds very bad" is added for row 14.
$ws.Range("G7").Value  = "Maximum, with an empty stage"
$ws.Range("G8").Value  = "Maximum, without animations."
$ws.Range("G9").Value  = "Achievable in very simple scenes"
$ws.Range("G10").Value = "Realistic target"
$ws.Range("G11").Value = "Realistic target"
$ws.Range("G12").Value = "Realistic target with some animations"
$ws.Range("G13").Value = "Minimum acceptable quality"
$ws.Range("G14").Value = "Sounds very bad"

# --- Wrap-text formatting follows the note text: row 9 no longer wraps
#     (short note moved up into it), every other note row wraps.
$ws.Range("G3").WrapText  = $true
$ws.Range("G4").WrapText  = $true
$ws.Range("G5").WrapText  = $true
$ws.Range("G6").WrapText  = $true
$ws.Range("G7").WrapText  = $true
$ws.Range("G8").WrapText  = $true
$ws.Range("G9").WrapText  = $false
$ws.Range("G10").WrapText = $true
$ws.Range("G11").WrapText = $true
$ws.Range("G12").WrapText = $true
$ws.Range("G13").WrapText = $true
$ws.Range("G14").WrapText = $true

# --- Row heights: the two "tall" (wrapped, two-line) rows move from
#     9 & 13 to 8 & 12, and a new tall row appears at 7.
$ws.Rows.Item(7).RowHeight  = 23.85
$ws.Rows.Item(8).RowHeight  = 23.85
$ws.Rows.Item(9).RowHeight  = 12.8
$ws.Rows.Item(12).RowHeight = 23.85
$ws.Rows.Item(13).RowHeight = 12.8

# --- Selection moves to G12.
$ws.Range("G12").Select() | Out-Null

$wb.Save()
